$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# Update the Kadjust value (H2) from 1 to 1.5
$ws.Range("H2").Value = 1.5

# Add explanatory comment to the Kadjust header cell (H1)
$commentText = "Kadjust for mesic peat. K_Sapric = K_mesic x 0.5" + [char]10 + "K_fibric = K_mesic x 2" + [char]10
$cell = $ws.Range("H1")
if ($cell.Comment -ne $null) {
    $cell.Comment.Delete()
}
$cell.AddComment($commentText)

$ws.Activate()
$ws.Range("H2").Select()
